$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "admin123"
$ws.Range("A3").Value = "admin2"
$ws.Range("A2").Value = "Admin1"
$ws.Range("B3").Value = "admin908"

$ws.Range("B3").Select()
